# RPA datasets push 2024-07-18
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("02_38커뮤니케이션(최근일자기준)")

# Insert a new row at row 3, pushing existing rows 3..21 down to 4..22
$ws.Rows.Item(3).Insert()

# Fill the new row 3 with the "이엔셀" entry
$ws.Cells.Item(3, 1).Value = "이엔셀"
$ws.Cells.Item(3, 2).Value = "2024.08.02~08.08"
$ws.Cells.Item(3, 3).Value = "13,600~15,300"
$ws.Cells.Item(3, 4).Value = "-"
$ws.Cells.Item(3, 5).Value = 21308
$ws.Cells.Item(3, 6).Value = "NH투자증권"

# The old "이엔셀" row has now shifted down to row 12 and is now a stale duplicate; remove it
$ws.Rows.Item(12).Delete()

# Update the 산일전기(유가) row's 확정공모가 (row 19, column D) from "-" to 35000
# (kept as text, matching the sibling cells in this column that store
# confirmed-offer-price as text rather than a number)
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "35000"
